$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.157.94'
$ws.Range('E2').Value = '  -2.19%  '
$ws.Range('D3').Value = '1.574.73'
$ws.Range('E3').Value = '  -1.77%  '
$ws.Range('E4').Value = '  -0.51%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '208.67'
$ws.Range('E5').Value = '  -1.53%  '
$ws.Range('E6').Value = '  -2.85%  '
$ws.Range('E7').Value = '  -0.46%  '
$ws.Range('E8').Value = '  -1.65%  '
$ws.Range('E9').Value = '  -1.29%  '
$ws.Range('E10').Value = '  -0.63%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0843'
$ws.Range('E11').Value = '  -0.45%  '
$ws.Range('D12').Value = '1.795.81'
$ws.Range('E12').Value = '  -1.80%  '
$ws.Range('D13').Value = '1.588.23'
$ws.Range('E13').Value = '  -0.94%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.05'
$ws.Range('E14').Value = '  -0.49%  '
$ws.Range('E15').Value = '  -2.09%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.36'
$ws.Range('E16').Value = '  -1.04%  '
$ws.Range('D17').Value = '26.135.06'
$ws.Range('E17').Value = '  -2.21%  '
$ws.Range('D18').Value = '0.0₃0726'
$ws.Range('E18').Value = '  -2.43%  '
$ws.Range('E19').Value = '  +1.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '207.72'
$ws.Range('E20').Value = '  -1.03%  '
$ws.Range('E22').Value = '  -1.23%  '
$ws.Range('E23').Value = '  -2.90%  '
$ws.Range('E24').Value = '  -2.70%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.55'
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('E26').Value = '  -0.29%  '
$ws.Range('E27').Value = '  -1.67%  '
$ws.Range('E28').Value = '  -1.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.20'
$ws.Range('E29').Value = '  -1.08%  '
$ws.Range('E30').Value = '  -0.55%  '
$ws.Range('E32').Value = '  -2.00%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.98'
$ws.Range('E33').Value = '  +0.38%  '
$ws.Range('D34').Value = '1.276.83'
$ws.Range('E34').Value = '  -1.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.616'
$ws.Range('E35').Value = '  +3.91%  '
$ws.Range('E37').Value = '  -0.90%  '
$ws.Range('E38').Value = '  -2.53%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.10'
$ws.Range('E39').Value = '  -10.65%  '
$ws.Range('E40').Value = '  -2.28%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  -0.47%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.56'
$ws.Range('E42').Value = '  +2.07%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.13'
$ws.Range('E43').Value = '  -2.72%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.763'
$ws.Range('E44').Value = '  -2.33%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '62.17'
$ws.Range('E45').Value = '  -1.41%  '
$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D46').Value = '1.708.74'
$ws.Range('E46').Value = '  -1.84%  '
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '88.90'
$ws.Range('E47').Value = '  -1.81%  '
$ws.Range('E48').Value = '  -3.11%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0103'
$ws.Range('E49').Value = '  +0.34%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.100'
$ws.Range('E50').Value = '  -1.95%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0506'
$ws.Range('E51').Value = '  -1.53%  '
